# Append a new data row (row 13) to the NIFTY options sheet, mirroring the
# structure/style of the existing rows (row 12 is the most recent one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 13

# --- Values that look like a date / percentage need an explicit text
# number format BEFORE the value is assigned, otherwise Excel's normal
# auto-detection will silently convert "2026-02-19" into a date serial
# and "100%" into the fraction 1 (with percent formatting). Every other
# cell in this table is plain text already, so it's unaffected.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2026-02-19"

$ws.Range("E$row").NumberFormat = "@"
$ws.Range("E$row").Value = "100%"

# Remaining text / string columns
$ws.Range("B$row").Value = "10:00:13"
$ws.Range("C$row").Value = "AVOID"
$ws.Range("D$row").Value = "AVOID"
$ws.Range("F$row").Value = "TRADEABLE"

# Numeric columns
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Value = 25773.7
$ws.Range("I$row").Value = 12.21
$ws.Range("J$row").Value = -1.12
$ws.Range("K$row").Value = 0
$ws.Range("L$row").Value = 46.5

$ws.Range("M$row").Value = "UNKNOWN"
$ws.Range("N$row").Value = 0
$ws.Range("O$row").Value = "UNKNOWN"
$ws.Range("P$row").Value = 0
$ws.Range("Q$row").Value = 0
$ws.Range("R$row").Value = 0
$ws.Range("S$row").Value = 0

$ws.Range("T$row").Value = "NONE"
$ws.Range("V$row").Value = 0
$ws.Range("W$row").Value = 0
$ws.Range("X$row").Value = 0
$ws.Range("Y$row").Value = 0
$ws.Range("Z$row").Value = 0
$ws.Range("AA$row").Value = 0
$ws.Range("AB$row").Value = 0

$ws.Range("AC$row").Value = "HARD VETO: CPR TRENDING DAY: Price 25773.70 within CPR (25736.60 - 25791.77) - SIDEWAYS/RANGE-BOUND BUT VERY NARROW CPR (0.214%) suggests trending day"
$ws.Range("AD$row").Value = "CPR TRENDING DAY: Price 25773.70 within CPR (25736.60 - 25791.77) - SIDEWAYS/RANGE-BOUND BUT VERY NARROW CPR (0.214%) suggests trending day"
$ws.Range("AE$row").Value = "Yes"

# U13 mirrors U12: present in the sheet as an explicit empty text value
# (not a formula, not a fully blank/unused cell). A lone "'" is Excel's
# text-prefix marker for an apostrophe-only entry and resolves to an
# empty string value once stored.
$ws.Range("U$row").Value = "'"

# Copy formatting from the row above so the new row matches the existing
# table styling (same fills/borders/number formats per column) -- this
# also normalises A13/E13 back from the temporary "@" text format to the
# shared style used by the rest of the column, while leaving their
# values as text.
$ws.Range("A12:AE12").Copy()
$ws.Range("A${row}:AE${row}").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
